$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force column D to be treated as text so numeric-looking price strings
# (e.g. "1.01") are not auto-converted to numbers, matching the original
# inlineStr text cells.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = '42.647.29'
$ws.Range("E2").Value = '  +0.48%  '
$ws.Range("D3").Value = '2.271.57'
$ws.Range("E3").Value = '  -0.77%  '
$ws.Range("E4").Value = '  -0.34%  '
$ws.Range("D5").Value = '308.54'
$ws.Range("E5").Value = '  -3.81%  '
$ws.Range("D6").Value = '102.66'
$ws.Range("E6").Value = '  +0.53%  '
$ws.Range("D7").Value = '0.619'
$ws.Range("E7").Value = '  -0.31%  '
$ws.Range("D8").Value = '1.01'
$ws.Range("E8").Value = '  +0.26%  '
$ws.Range("D9").Value = '0.596'
$ws.Range("E9").Value = '  -1.65%  '
$ws.Range("D10").Value = '38.37'
$ws.Range("E10").Value = '  -2.52%  '
$ws.Range("D11").Value = '0.0893'
$ws.Range("E11").Value = '  -1.08%  '
$ws.Range("D12").Value = '8.14'
$ws.Range("E12").Value = '  -1.82%  '
$ws.Range("E13").Value = '  +0.96%  '
$ws.Range("D14").Value = '0.967'
$ws.Range("E14").Value = '  +0.60%  '
$ws.Range("D15").Value = '14.93'
$ws.Range("E15").Value = '  -1.32%  '
$ws.Range("D16").Value = '2.613.62'
$ws.Range("E16").Value = '  -0.79%  '
$ws.Range("D17").Value = '2.273.92'
$ws.Range("E17").Value = '  -0.82%  '
$ws.Range("D18").Value = '42.512.77'
$ws.Range("E18").Value = '  +0.35%  '
$ws.Range("D19").Value = '7.17'
$ws.Range("E19").Value = '  -2.79%  '
$ws.Range("D20").Value = '0.0000103'
$ws.Range("E20").Value = '  -1.78%  '
$ws.Range("D21").Value = '13.02'
$ws.Range("E21").Value = '  +2.35%  '
$ws.Range("D22").Value = '72.50'
$ws.Range("E22").Value = '  -0.63%  '
$ws.Range("D23").Value = '3.35'
$ws.Range("E23").Value = '  -5.95%  '
$ws.Range("D24").Value = '260.93'
$ws.Range("E24").Value = '  -2.72%  '
$ws.Range("D25").Value = '2.15'
$ws.Range("E25").Value = '  -2.99%  '
$ws.Range("E26").Value = '  +0.21%  '
$ws.Range("D27").Value = '10.55'
$ws.Range("E27").Value = '  -2.80%  '
$ws.Range("E28").Value = '  -0.03%  '
$ws.Range("D29").Value = '6.81'
$ws.Range("E29").Value = '  +12.75%  '
$ws.Range("D30").Value = '22.01'
$ws.Range("E30").Value = '  -2.17%  '
$ws.Range("D31").Value = '35.35'
$ws.Range("E31").Value = '  -6.66%  '
$ws.Range("D32").Value = '163.58'
$ws.Range("E32").Value = '  -0.83%  '
$ws.Range("D33").Value = '0.0845'
$ws.Range("E33").Value = '  -3.01%  '
$ws.Range("E34").Value = '  -2.72%  '
$ws.Range("E35").Value = '  +2.01%  '
$ws.Range("D36").Value = '0.110'
$ws.Range("E36").Value = '  -2.80%  '
$ws.Range("D37").Value = '4.47'
$ws.Range("E37").Value = '  -2.26%  '
$ws.Range("D38").Value = '0.0344'
$ws.Range("E38").Value = '  -2.95%  '
$ws.Range("D39").Value = '3.65'
$ws.Range("E39").Value = '  -0.55%  '
$ws.Range("D40").Value = '2.71'
$ws.Range("E40").Value = '  -1.84%  '
$ws.Range("D41").Value = '1.53'
$ws.Range("E41").Value = '  +0.06%  '
$ws.Range("D42").Value = '97.83'
$ws.Range("E42").Value = '  +7.60%  '
$ws.Range("D43").Value = '1.01'
$ws.Range("E43").Value = '  +0.25%  '
$ws.Range("D44").Value = '67.77'
$ws.Range("E44").Value = '  -0.65%  '
$ws.Range("D45").Value = '0.223'
$ws.Range("E45").Value = '  -0.48%  '
$ws.Range("B46").Value = 'Maker'
$ws.Range("C46").Value = 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr'
$ws.Range("D46").Value = '1.713.98'
$ws.Range("E46").Value = '  +6.59%  '
$ws.Range("B47").Value = 'Celestia'
$ws.Range("C47").Value = 'https://coinranking.com/coin/YQcD0lBl7+celestia-tia'
$ws.Range("D47").Value = '11.75'
$ws.Range("E47").Value = '  -4.00%  '
$ws.Range("D48").Value = '108.89'
$ws.Range("E48").Value = '  -4.59%  '
$ws.Range("D49").Value = '75.24'
$ws.Range("E49").Value = '  -5.53%  '
$ws.Range("B50").Value = 'FraxShare'
$ws.Range("C50").Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range("D50").Value = '8.56'
$ws.Range("E50").Value = '  -4.44%  '
$ws.Range("B51").Value = 'THORChain'
$ws.Range("C51").Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range("D51").Value = '5.09'
$ws.Range("E51").Value = '  -2.57%  '
